# Applies the two changes described in the diff:
#  1. Merge the two runs " instances" and " when there is extra whitespace"
#     into a single run " instances when there is extra whitespace".
#  2. Remove the standalone paragraph holding the "line: 3" test-data row
#     (the "line: " / "3" runs) from the candidatesLine table cell.

$d = $word.ActiveDocument

# --- Change 1: merge the split sentence back into one run -------------------
$d.Content.Find.Execute(
    " instances when there is extra whitespace",  # FindText
    $false,                                        # MatchCase
    $false,                                        # MatchWholeWord
    $false,                                        # MatchWildcards
    $false,                                        # MatchSoundsLike
    $false,                                        # MatchAllWordForms
    $true,                                          # Forward
    1,                                              # Wrap (wdFindContinue)
    $false,                                         # Format
    " instances when there is extra whitespace",   # ReplaceWith
    2                                               # Replace (wdReplaceAll)
) | Out-Null

# --- Change 2: delete the "line: 3" paragraph --------------------------------
# Note: this paragraph is the last one in its table cell, so Range.Text
# carries a trailing cell-mark (chr(7)) after the paragraph mark - trim
# trailing control characters before comparing.
$paras = $d.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text.TrimEnd("`r", "`a")
    if ($txt -eq "line: 3") {
        $p.Range.Delete()
        break
    }
}
